$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 136; this shifts existing rows 136..199 down to 137..200.
$ws.Rows(136).Insert()

# Copy the content of what is now row 137 (the old row 136) into the newly
# inserted (blank) row 136, then update the date to the new reporting date.
$ws.Range("A136:T136").Value2 = $ws.Range("A137:T137").Value2

$ws.Cells.Item(136, 4).Value2 = 44460
